$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Capitalize the first letter of each user's name in column A (rows 2-21)
$ws.Range("A2").Value  = "Soumya"
$ws.Range("A3").Value  = "Kiran"
$ws.Range("A4").Value  = "Pratyusha"
$ws.Range("A5").Value  = "Shubha"
$ws.Range("A6").Value  = "Rajni"
$ws.Range("A7").Value  = "Akansha"
$ws.Range("A8").Value  = "Hindu"
$ws.Range("A9").Value  = "Bindu"
$ws.Range("A10").Value = "Sindhu"
$ws.Range("A11").Value = "Nisha"
$ws.Range("A12").Value = "Pallavi"
$ws.Range("A13").Value = "Harika"
$ws.Range("A14").Value = "Akhila"
$ws.Range("A15").Value = "Monika"
$ws.Range("A16").Value = "Priyanka"
$ws.Range("A17").Value = "Shradda"
$ws.Range("A18").Value = "Nithya"
$ws.Range("A19").Value = "Poornima"
$ws.Range("A20").Value = "Chethana"
$ws.Range("A21").Value = "Akshar"

# 2. Refresh the cached "display" text of the password-column hyperlinks so it
#    matches the (already updated) "@1234" cell text instead of the stale
#    "@123" text. The hyperlink targets themselves stay exactly as they were.
#    The engine does not allow editing hyperlinks that were loaded from the
#    workbook in-place, so every hyperlink on the sheet is removed and then
#    re-created in its original order (this keeps the relationship ids and
#    ordering stable) with the corrected display text.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"),  "mailto:Soumya@123",      "", "", "Soumya@1234")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"),  "mailto:Kiran@123",       "", "", "Kiran@1234")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"),  "mailto:Pratyusha@123",   "", "", "Pratyusha@1234") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"),  "mailto:Shubha@123",      "", "", "Shubha@1234")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"),  "mailto:Rajni@123",       "", "", "Rajni@1234")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("E7"),  "mailto:Akansha@123",     "", "", "Akansha@1234")   | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"),  "mailto:Hindu@123",       "", "", "Hindu@1234")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"),  "mailto:Bindu@123",       "", "", "Bindu@1234")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), "mailto:Sindhu@123",      "", "", "Sindhu@1234")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E11"), "mailto:Nishi@123",       "", "", "Nishi@1234")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("E12"), "mailto:Pallavi@123",     "", "", "Pallavi@1234")   | Out-Null
$ws.Hyperlinks.Add($ws.Range("E13"), "mailto:Harika@123",      "", "", "Harika@1234")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E14"), "mailto:Akhil@123",       "", "", "Akhil@1234")     | Out-Null
$ws.Hyperlinks.Add($ws.Range("E15"), "mailto:Monika@123",      "", "", "Monika@1234")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E16"), "mailto:Priyanka@123",    "", "", "Priyanka@1234")  | Out-Null

$ws.Hyperlinks.Add($ws.Range("B17"), "mailto:shradda@mail.com",   "", "", "shradda@mail.com")   | Out-Null
$ws.Hyperlinks.Add($ws.Range("E17"), "mailto:Shradda@123",        "", "", "Shradda@1234")        | Out-Null
$ws.Hyperlinks.Add($ws.Range("B18"), "mailto:nithya@gmail.com",   "", "", "nithya@gmail.com")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E18"), "mailto:Nithya@123",         "", "", "Nithya@1234")         | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:poornima@gmail.com", "", "", "poornima@gmail.com")  | Out-Null
$ws.Hyperlinks.Add($ws.Range("E19"), "mailto:Poornima@123",       "", "", "Poornima@1234")       | Out-Null
$ws.Hyperlinks.Add($ws.Range("B20"), "mailto:chethana@gmail.com", "", "", "chethana@gmail.com")  | Out-Null
$ws.Hyperlinks.Add($ws.Range("E20"), "mailto:Chethana@123",       "", "", "Chethana@1234")       | Out-Null
$ws.Hyperlinks.Add($ws.Range("B21"), "mailto:akshar@gmail.com",   "", "", "akshar@gmail.com")    | Out-Null
$ws.Hyperlinks.Add($ws.Range("E21"), "mailto:Akshar@123",         "", "", "Akshar@123")          | Out-Null

# 3. Move the active selection to A21 (was I8)
$ws.Range("A21").Select() | Out-Null
